$d = $word.ActiveDocument

function Set-ParagraphXml($searchText, $xml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $rng.Expand(4) | Out-Null
    $rng.InsertXML($xml)
}

# --- Names header: split "Ana Luisa Giaquinto Zolio" into spell-checked runs
# and insert a new "LINK GIT:" paragraph after the "Julia Martins" line ---
$rngName = $d.Content
$foundName = $rngName.Find.Execute("Ana Luisa Giaquinto", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundName) { throw "Could not find Ana Luisa paragraph" }
$rngName.Expand(4) | Out-Null
$rngJulia = $d.Content
$foundJulia = $rngJulia.Find.Execute("Julia Martins de Almeida Antunes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundJulia) { throw "Could not find Julia Martins paragraph" }
$rngJulia.Expand(4) | Out-Null
$comboRng = $d.Range($rngName.Start, $rngJulia.End)
$comboXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Ana </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Luisa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Giaquinto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Zólio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> RM99348</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Julia Martins de Almeida Antunes RM98601</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">LINK GIT: </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>https://github.com/JuMartinsDev/GS_Governan-a</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$comboRng.InsertXML($comboXml)

# --- a1 ---
$xml_a1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>a.1. Crie um repositório público no GITHUB.com para o seu projeto (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ɵre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> print </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> da tela do repositório criado com a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Branch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> e cole no documento de resposta da sua prova);</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Set-ParagraphXml "a.1. Crie um repositório" $xml_a1

# --- a2 ---
$xml_a2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve"> a.2. Usando o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> CMD ou BASH, clone o repositório em uma pasta local projeto (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ɵre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> print </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> da tela que mostra a execução do comando de clonagem e cole no documento de resposta da sua prova); </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Set-ParagraphXml "a.2. Usando o Git CMD" $xml_a2

# --- a3 ---
$xml_a3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>a.3. Inicialize o GIT FLOW projeto (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ɵre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> print </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> da tela que mostra a execução do comando e cole no documento de resposta da sua prova).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Set-ParagraphXml "a.3. Inicialize o GIT FLOW" $xml_a3

# --- termine ---
$xml_termine = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Termine a edição do arquivo de resposta da sua prova, colocando o nome e RM dos integrantes do seu grupo no alto da 1ª página (se </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ɵver</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> feito em grupo) e... </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Set-ParagraphXml "Termine a edição do arquivo" $xml_termine

# --- c1 ---
$xml_c1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>c.1. Usando o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> CMD ou BASH, faça o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>commit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> do documento na pasta GIT local (clonada) na </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Branch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>develop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ɵre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> print </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> da tela que mostra a execução do comando e cole no documento de resposta da sua prova); </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Set-ParagraphXml "c.1. Usando o Git CMD" $xml_c1

# --- c2 ---
$xml_c2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>c.2. Faça a atualização da </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Branch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> após atualizar a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>develop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ɵre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> print </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> da tela que mostra a execução do comando e cole no documento de resposta da sua prova); </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Set-ParagraphXml "c.2. Faça a atualização" $xml_c2

# --- c3 ---
$xml_c3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>c.3. Atualize a cópia remota, enviando as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Branches</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>develop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para o GITHUB (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ɵre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> print </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>screen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> da tela que mostra a execução do comando de sincronização de dados e da tela do GITHUB com as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>branches</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> atualizadas).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Set-ParagraphXml "c.3. Atualize a cópia remota" $xml_c3

Write-Host "All edits applied successfully"